$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3130
$ws1.Range("F4").Value = 1090
$ws1.Range("F7").Value = 271
$ws1.Range("F8").Value = 35
$ws1.Range("F9").Value = 1124
$ws1.Range("F10").Value = 15630
$ws1.Range("F11").Value = 236
$ws1.Range("F12").Value = 172
$ws1.Range("F14").Value = 6159
$ws1.Range("F15").Value = 620
$ws1.Range("F18").Value = 6
$ws1.Range("F20").Value = 1262
$ws1.Range("F23").Value = 12
$ws1.Range("F24").Value = 12
$ws1.Range("F28").Value = 25
$ws1.Range("F29").Value = 4994
$ws1.Range("F30").Value = 482
$ws1.Range("F31").Value = 11033
$ws1.Range("F32").Value = 1226
$ws1.Range("F34").Value = 115
$ws1.Range("F35").Value = 163
$ws1.Range("F36").Value = 3797
$ws1.Range("F37").Value = 263

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3130
$ws4.Range("F5").Value = 1090
$ws4.Range("F8").Value = 271
$ws4.Range("F9").Value = 35
$ws4.Range("F10").Value = 1124
$ws4.Range("F11").Value = 15630
$ws4.Range("F12").Value = 236
$ws4.Range("F13").Value = 172
$ws4.Range("F15").Value = 6159
$ws4.Range("F16").Value = 620
$ws4.Range("F19").Value = 6
$ws4.Range("F21").Value = 1262
$ws4.Range("F24").Value = 12
$ws4.Range("F25").Value = 12
$ws4.Range("F29").Value = 25
$ws4.Range("F30").Value = 4994
$ws4.Range("F31").Value = 482
$ws4.Range("F33").Value = 11033
$ws4.Range("F34").Value = 1226
$ws4.Range("F36").Value = 115
$ws4.Range("F37").Value = 163
$ws4.Range("F38").Value = 3797
$ws4.Range("F39").Value = 263
